$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("E1").Value = "Notes"

# New data row (row 3) - copy formatting from row 2 so styles are reused
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("C2").Copy($ws.Range("C3"))

$ws.Range("A3").Value = 42369
$ws.Range("B3").Value = 0.625
$ws.Range("C3").Value = 0.66666666666666663

$ws.Range("E3").Value = "Learning about GIT"

$ws.Range("C4").Select()
